$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the three changed values in row 3
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 13

# Leave the selection on the last-edited cell (E3), matching the saved view state
$ws.Range("E3").Select()
